# approach section; README update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recording Sheet")

$url = "https://github.com/alex-kalinka-cruk/fgc_crispr_pipeline_UAT/blob/master/analysis/UAT-analysis.Rmd"

# Column B ("Test Case/Script") now links out to the UAT-analysis.Rmd README
# on GitHub instead of bare text. Apply the Hyperlink style first so the
# generated style matches a plain "Hyperlink" cell style (no extra wrap
# carried over from the old formatting), then attach the actual hyperlinks.
$b2 = $ws.Cells.Item(2, 2)
$b2.Style = "Hyperlink"
$b2.Value = $url
[void]$ws.Hyperlinks.Add($b2, $url)

# New "Test Step Header" (column D) content describing the approach taken
# for each test case.
$ws.Cells.Item(2, 4).Value = "Manual test of data transformation for multi-SLX use case."
$ws.Cells.Item(3, 4).Value = "# Concordance of pipeline outputs with v1"
$ws.Cells.Item(4, 4).Value = "# Concordance of pipeline outputs with v2"
$ws.Cells.Item(5, 4).Value = "# Concordance of pipeline outputs with v3"
$ws.Cells.Item(6, 4).Value = "# Concordance of pipeline outputs with v4"
$ws.Cells.Item(7, 4).Value = "# Concordance of pipeline outputs with v5"
$ws.Cells.Item(8, 4).Value = "# Concordance of pipeline outputs with v6"
$ws.Cells.Item(9, 4).Value = "# Concordance of pipeline outputs with v7"

$bRest = $ws.Range($ws.Cells.Item(3, 2), $ws.Cells.Item(9, 2))
$bRest.Style = "Hyperlink"
$bRest.Value = $url
[void]$ws.Hyperlinks.Add($bRest, $url, "", "", $url)

# Restore the selected cell recorded in the sheet view.
[void]$ws.Range("A14").Select()
